# This script reproduces the commit that appended two new daily price
# records into the "Berenjena" (eggplant) price sheet:
#
#   - A new record (2022-01-06, serial 44567) is inserted at row 129,
#     pushing the old rows 129..200 down to 130..201.
#   - A new record (2022-01-07, serial 44568) is inserted at row 202,
#     pushing the old rows 201..210 (which by then sit at 202..211) down
#     to 203..212.
#
# The workbook's used range therefore grows from A1:R210 to A1:R212.
# Every other row keeps its original data; it is simply relocated by the
# native row-insert operation below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-Row($RowIndex, $Fecha, $Volumen, $PrecioMinimo, $PrecioMaximo, $PrecioPromedio, $Origen, $PrecioKg) {
    # Insert a blank row here; Excel shifts the current row (and below) down.
    $ws.Rows($RowIndex).Insert()

    $ws.Cells.Item($RowIndex, 1).Value = 3
    $ws.Cells.Item($RowIndex, 2).Value = "Femacal de La Calera"
    $ws.Cells.Item($RowIndex, 3).Value = "Coquimbo"

    $dCell = $ws.Cells.Item($RowIndex, 4)
    $dCell.Value = $Fecha
    $dCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($RowIndex, 5).Value = 5
    $ws.Cells.Item($RowIndex, 6).Value = 100112001
    $ws.Cells.Item($RowIndex, 7).Value = "Berenjena"
    $ws.Cells.Item($RowIndex, 8).Value = "Sin especificar"
    $ws.Cells.Item($RowIndex, 9).Value = "Primera"
    $ws.Cells.Item($RowIndex, 10).Value = $Volumen
    $ws.Cells.Item($RowIndex, 11).Value = $PrecioMinimo
    $ws.Cells.Item($RowIndex, 12).Value = $PrecioMaximo
    $ws.Cells.Item($RowIndex, 13).Value = $PrecioPromedio
    $ws.Cells.Item($RowIndex, 14).Value = "`$/caja 60 unidades"
    $ws.Cells.Item($RowIndex, 15).Value = $Origen
    $ws.Cells.Item($RowIndex, 16).Value = $PrecioKg
    $ws.Cells.Item($RowIndex, 17).Value = 60
    $ws.Cells.Item($RowIndex, 18).Value = "Hortaliza"
}

# New record #1 -> lands at row 129 (2022-01-06).
Add-Row 129 44567 110 7500 8000 7773 "Región de Arica y Parinacota" 130

# New record #2 -> lands at row 202 (2022-01-07).
Add-Row 202 44568 120 7500 8000 7750 "Región de Arica y Parinacota" 129
